$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Order id header + value saved for the Android app to read back.
$ws.Range("A1").Value = "oid"

# The order id is stored as text (it is display/lookup data, not used in
# math), so force a text format before assigning it - this also preserves
# the leading space in the original value.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = " 188285"

$ws.Range("B3").Select()
